$d = $word.ActiveDocument

$replacements = @(
    @("2025-08-14 Thursday", "2025-08-15 Friday"),
    @("618×6=3708", "652×6=3912"),
    @("643×9=5787", "454×7=3178"),
    @("670×8=5360", "435×3=1305"),
    @("440×5=2200", "487×7=3409"),
    @("778×3=2334", "990×6=5940"),
    @("486×7=3402", "442×4=1768"),
    @("390×3=1170", "579×7=4053"),
    @("541×8=4328", "819×3=2457"),
    @("266×6=1596", "683×2=1366"),
    @("779×5=3895", "892×5=4460"),
    @("990×3=2970", "372×7=2604"),
    @("170×4=680",  "988×6=5928"),
    @("740×8=5920", "581×6=3486"),
    @("634×9=5706", "153×5=765"),
    @("162×7=1134", "837×4=3348"),
    @("243×5=1215", "799×4=3196"),
    @("249×8=1992", "502×3=1506"),
    @("471×3=1413", "474×4=1896"),
    @("755×5=3775", "725×3=2175"),
    @("966×6=5796", "198×6=1188"),
    @("722×9=6498", "980×8=7840"),
    @("695×8=5560", "177×6=1062"),
    @("481×2=962",  "994×8=7952"),
    @("829×5=4145", "999×2=1998"),
    @("692×3=2076", "605×7=4235")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
